$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A60").Value = "NIRampTime"
$ws.Range("B60").Value = "ms"

$wb.Application.CalculateFullRebuild()
